$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 434
    $ws.Range("F3").Value = 1456
    $ws.Range("F6").Value = 2136
    $ws.Range("F8").Value = 1297
    $ws.Range("F10").Value = 122
}
